$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix capitalization of search input text
$ws.Range("A2").Value = "Vladimir putine"

# 2. Update match rate for row 5 (2nd result)
$ws.Range("B5").Value = "2. (95.84%) -  Vladimir Vladimirovich Putin"

# 3. Update row 6 (3rd result) to reflect new match rate/source list, replacing
#    the old "FR - DGT" source with "GB - HMT" (same as rows 3 & 4) and a new link id
$ws.Range("B6").Value = "3. (86.75%) -  Vladimir Vladimirovich Putin"
$ws.Range("C6").Value = "GB - Liste consolidée des sanctions financières du Royaume-Uni (HMT)"
$ws.Range("G6").Value = "todoByFrontDev/655190123456789012345678"

# 4. Remove now-obsolete result rows 7-13 (results dropped from 10 to 4 total)
$ws.Rows("7:13").Delete()

Write-Output "edit complete"
